# Applies the updated simulation-result values to the workbook.
# The workbook holds static result data (no formulas), so the edit is a
# direct set of cell-value writes on the affected sheets.

$wb = $excel.ActiveWorkbook

# --- pcroprep (sheet1) ---------------------------------------------------
$ws = $wb.Worksheets.Item("pcroprep")
$ws.Range("D35").Value = 1.4105556681109319
$ws.Range("F35").Value = 5.1292933385852058
$ws.Range("G35").Value = -236.27070666141481

$ws.Range("D39").Value = 1222.3405183637792
$ws.Range("F39").Value = 384.03026368106066
$ws.Range("G39").Value = -422.56973631893942

# --- pdietrep (sheet4) ----------------------------------------------------
$ws = $wb.Worksheets.Item("pdietrep")
$ws.Range("E6").Value = 1553.2124943666656
$ws.Range("F6").Value = -604.57301983901721
$ws.Range("G6").Value = 71.981783367306974

$ws.Range("E7").Value = 59.035912064222963
$ws.Range("F7").Value = -11.092278812127596
$ws.Range("G7").Value = 84.182853324014303

$ws.Range("E8").Value = 24.472570235374661
$ws.Range("F8").Value = -40.260995190795811
$ws.Range("G8").Value = 37.805070791730088

$ws.Range("E9").Value = 257.41756835777619
$ws.Range("F9").Value = -66.250258773076155
$ws.Range("G9").Value = 79.531404353546534

# --- pradar (sheet5) ------------------------------------------------------
$ws = $wb.Worksheets.Item("pradar")
$ws.Range("D15").Value = 5.1292933385852058
$ws.Range("E15").Value = 2.124810827914335
$ws.Range("F15").Value = -236.27070666141481

# --- plandrep (sheet6) -----------------------------------------------------
$ws = $wb.Worksheets.Item("plandrep")
$ws.Range("S11").Value = 10.434000000000061

# --- plaborrep (sheet7) ----------------------------------------------------
$ws = $wb.Worksheets.Item("plaborrep")
$ws.Range("R3").Value = 0.0085037100000000341
$ws.Range("AF3").Value = 1.3863642563050813

# --- pfertrep (sheet8) -----------------------------------------------------
$ws = $wb.Worksheets.Item("pfertrep")
$ws.Range("S5").Value = 2869.3500000000167
$ws.Range("Z5").Value = 414894.23386399995

$ws.Range("S6").Value = 1721.6100000000101
$ws.Range("Z6").Value = 518933.48162400001

$ws.Range("S7").Value = 2744.1420000000162
$ws.Range("Z7").Value = 478309.46967199992
